$d = $word.ActiveDocument

$pairs = @(
    @("510÷8=", "164÷7="),
    @("103÷4=", "890÷6="),
    @("737÷2=", "492÷6="),
    @("171÷2=", "565÷9="),
    @("403÷2=", "423÷6="),
    @("478÷8=", "340÷4="),
    @("334÷5=", "922÷7="),
    @("100÷9=", "996÷4="),
    @("286÷8=", "100÷9="),
    @("838÷6=", "225÷4="),
    @("925÷7=", "188÷8="),
    @("146÷3=", "394÷5="),
    @("231÷7=", "485÷9="),
    @("275÷8=", "651÷4="),
    @("758÷2=", "275÷3="),
    @("372÷2=", "465÷6="),
    @("436÷5=", "234÷2="),
    @("580÷4=", "571÷9="),
    @("928÷8=", "624÷9="),
    @("894÷5=", "731÷2="),
    @("240÷2=", "641÷5="),
    @("279÷3=", "458÷2="),
    @("352÷2=", "793÷5="),
    @("192÷6=", "332÷6="),
    @("224÷3=", "363÷2="),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
